$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.538860678672791
$ws.Range("B1").Value = 3.791602611541748
$ws.Range("C1").Value = 3.311583757400513
$ws.Range("D1").Value = 1.474736571311951
$ws.Range("E1").Value = 1.012883305549622
